# Update gh-pages to output generated at 456a3b4
# Applies the refreshed "想去人数" (want-to-go count) figures, a refreshed
# cover-image URL, and a corrected time range that came out of the new
# data pull for 江西-漫展信息.xlsx.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" --------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5133
$ws1.Range("F5").Value = 7435
$ws1.Range("F11").Value = 27
$ws1.Range("F12").Value = 4313
$ws1.Range("F13").Value = 1758
$ws1.Range("F16").Value = 2916
$ws1.Range("F19").Value = 206
$ws1.Range("F20").Value = 499
$ws1.Range("F21").Value = 434
$ws1.Range("F22").Value = 456
$ws1.Range("F24").Value = 99
$ws1.Range("F25").Value = 1692
$ws1.Range("I30").Value = "//i1.hdslb.com/bfs/openplatform/202407/wzeWuUHS1721735596448.jpeg"
$ws1.Range("F32").Value = 515
$ws1.Range("F34").Value = 61
$ws1.Range("F35").Value = 105
$ws1.Range("F36").Value = 63
$ws1.Range("F37").Value = 2879
$ws1.Range("F38").Value = 703
$ws1.Range("F40").Value = 66
$ws1.Range("F42").Value = 24

# ---- Sheet "演出" --------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("E3").Value = "2024.08.17 13:30-08.17 15:30"

# ---- Sheet "全部类型" -----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5133
$ws4.Range("F5").Value = 7435
$ws4.Range("F11").Value = 27
$ws4.Range("F12").Value = 4313
$ws4.Range("F13").Value = 1758
$ws4.Range("F16").Value = 2916
$ws4.Range("F19").Value = 206
$ws4.Range("F20").Value = 499
$ws4.Range("F21").Value = 434
$ws4.Range("F22").Value = 456
$ws4.Range("F25").Value = 99
$ws4.Range("F26").Value = 1692
$ws4.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202407/wzeWuUHS1721735596448.jpeg"
$ws4.Range("F33").Value = 515
$ws4.Range("F35").Value = 61
$ws4.Range("F36").Value = 105
$ws4.Range("F37").Value = 63
$ws4.Range("F38").Value = 2879
$ws4.Range("E39").Value = "2024.08.17 13:30-08.17 15:30"
$ws4.Range("F40").Value = 703
$ws4.Range("F42").Value = 66
$ws4.Range("F44").Value = 24
